$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused B2 value (row 2 "Buying Opportunity" column)
$ws.Range("B2").Value = ""

# Update "support Zone" column (C) with new tickers
$ws.Range("C2").Value  = "NSE:BOROLTD"
$ws.Range("C3").Value  = "NSE:CANTABIL"
$ws.Range("C4").Value  = "NSE:CONCORDBIO"
$ws.Range("C5").Value  = "NSE:DCMNVL"
$ws.Range("C6").Value  = "NSE:ERIS"
$ws.Range("C7").Value  = "NSE:EXCELINDUS"
$ws.Range("C8").Value  = "NSE:GINNIFILA"
$ws.Range("C9").Value  = "NSE:ORCHPHARMA"
$ws.Range("C10").Value = "NSE:PPAP"

# Update "Short buildup" column (E) with new tickers
$ws.Range("E2").Value  = "NSE:DEEPAKNTR"
$ws.Range("E3").Value  = "NSE:HAL"
$ws.Range("E4").Value  = "NSE:INDHOTEL"
$ws.Range("E5").Value  = "NSE:JIOFIN"
$ws.Range("E6").Value  = "NSE:KALYANKJIL"
$ws.Range("E7").Value  = "NSE:LAURUSLABS"
$ws.Range("E8").Value  = "NSE:MARUTI"
$ws.Range("E9").Value  = "NSE:MUTHOOTFIN"

# Row 10's Short buildup entry is removed (no longer populated)
$ws.Range("E10").Value = ""

# Remove rows 11-19 entirely, shrinking the used range to A1:F10
$ws.Range("A11:F19").EntireRow.Delete()
